$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text so that decimal-looking
# values (e.g. "326.03", "1.001") are not silently re-interpreted as numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Updated cryptocurrency price/volume data (includes TRON/ShibaInu row swap)
$ws.Range("D2").Value = "29.508.58"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "1.919.18"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "326.03"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "0.4743"
$ws.Range("E7").Value = "  +2.58%  "
$ws.Range("D8").Value = "0.4097"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "47.82"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "0.08051"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "22.52"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("D13").Value = "1.927.12"
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D15").Value = "7.156"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "89.57"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.06596"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.00001030"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "17.79"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "29.519.18"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").Value = "5.548"
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").Value = "2.205"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "2.103.28"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").Value = "154.76"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "19.88"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").Value = "6.028"
$ws.Range("E29").Value = "  +11.28%  "
$ws.Range("D30").Value = "2.131"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "117.77"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "1.063"
$ws.Range("E32").Value = "  +8.49%  "
$ws.Range("D33").Value = "0.09553"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "1.434"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "3.563"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").Value = "5.408"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").Value = "0.06120"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "8.321"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "1.174"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("D42").Value = "2.553"
$ws.Range("E42").Value = "  +10.99%  "
$ws.Range("D43").Value = "0.1847"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").Value = "0.08017"
$ws.Range("E45").Value = "  +14.43%  "
$ws.Range("D46").Value = "1.285"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").Value = "0.5560"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").Value = "12.16"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D49").Value = "1.940"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("D50").Value = "113.18"
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("D51").Value = "44.89"
$ws.Range("E51").Value = "  -0.95%  "

# Restore the default (Normal) cell style now that the text values are locked in,
# so no stray number-format styling is left behind on the data cells.
$dataRange.Style = "Normal"
